# Updated cryptos list on Tue Jul 30 18:31:12 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep Price/Volume columns as text so values like "0.630" or "436.90" do not
# get silently coerced into numbers (which would drop trailing zeros).
$ws.Columns("D").NumberFormat = "@"
$ws.Columns("E").NumberFormat = "@"

function Set-Row($Row, $B, $C, $D, $E) {
    if ($B -ne $null) { $ws.Range("B$Row").Value = $B }
    if ($C -ne $null) { $ws.Range("C$Row").Value = $C }
    if ($D -ne $null) { $ws.Range("D$Row").Value = $D }
    if ($E -ne $null) { $ws.Range("E$Row").Value = $E }
}

# Row 2 - Bitcoin
Set-Row 2 $null $null "65.943.84" "  -2.19%  "
# Row 3 - Ethereum
Set-Row 3 $null $null "3.290.42" "  -0.66%  "
# Row 4 - TetherUSD
Set-Row 4 $null $null $null "  +0.00%  "
# Row 5 - BNB
Set-Row 5 $null $null "573.95" "  -0.65%  "
# Row 6 - Solana
Set-Row 6 $null $null "178.07" "  -4.48%  "
# Row 7 - XRP
Set-Row 7 $null $null "0.630" "  +4.51%  "
# Row 8 - USDC
Set-Row 8 $null $null $null "  +0.02%  "
# Row 9 - Dogecoin
Set-Row 9 $null $null $null "  -2.41%  "
# Row 10 - Toncoin
Set-Row 10 $null $null "6.71" "  +0.90%  "
# Row 11 - Cardano
Set-Row 11 $null $null "0.401" "  -2.26%  "
# Row 12 - WrappedliquidstakedEther2.0
Set-Row 12 $null $null "3.866.25" "  -0.67%  "
# Row 13 - TRON
Set-Row 13 $null $null $null "  -3.56%  "
# Row 14 - Avalanche
Set-Row 14 $null $null "26.66" "  -2.89%  "
# Row 15 - WrappedBTC
Set-Row 15 $null $null "66.022.66" "  -2.32%  "

# Row 16 & 17 swap: WrappedEther <-> ShibaInu
Set-Row 16 "ShibaInu" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib" "0.0000163" "  -2.07%  "
Set-Row 17 "WrappedEther" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" "3.285.97" "  -0.75%  "

# Row 18 - BitcoinCash
Set-Row 18 $null $null "436.90" "  -1.65%  "
# Row 19 - Polkadot
Set-Row 19 $null $null "5.59" "  -2.13%  "
# Row 20 - Chainlink
Set-Row 20 $null $null "13.29" "  -1.76%  "
# Row 21 - Uniswap
Set-Row 21 $null $null "7.41" "  -4.25%  "
# Row 22 - Litecoin
Set-Row 22 $null $null "72.57" "  -1.73%  "
# Row 23 - Dai
Set-Row 23 $null $null "0.999" "  -0.04%  "

# Row 24 & 25 swap: WrappedeETH <-> Polygon
Set-Row 24 "Polygon" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic" "0.512" "  -0.25%  "
Set-Row 25 "WrappedeETH" "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth" "3.435.29" "  -0.60%  "

# Row 26 - PEPE
Set-Row 26 $null $null $null "  -4.88%  "
# Row 27 - Kaspa
Set-Row 27 $null $null "0.194" "  +2.93%  "
# Row 28 - InternetComputer(DFINITY)
Set-Row 28 $null $null "8.94" "  -1.20%  "
# Row 29 - Binance-PegBSC-USD
Set-Row 29 $null $null "0.999" "  +0.12%  "
# Row 30 - PancakeSwap
Set-Row 30 $null $null "1.94" "  -1.96%  "
# Row 31 - EthereumClassic
Set-Row 31 $null $null "22.36" "  -2.19%  "
# Row 32 - USDe
Set-Row 32 $null $null $null "  +0.09%  "
# Row 33 - NEARProtocol
Set-Row 33 $null $null "5.15" "  -3.41%  "
# Row 34 - Aptos
Set-Row 34 $null $null $null "  -2.45%  "
# Row 35 - Fetch.AI
Set-Row 35 $null $null $null "  -3.89%  "

# Row 36 & 37 swap: ImmutableX <-> Monero
Set-Row 36 "Monero" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" "158.39" "  -2.85%  "
Set-Row 37 "ImmutableX" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" "1.45" "  -4.37%  "

# Row 38 - EnergySwap
Set-Row 38 $null $null "26.79" "  -1.83%  "
# Row 39 - Stacks
Set-Row 39 $null $null "1.79" "  -3.58%  "
# Row 40 - Maker
Set-Row 40 $null $null "2.785.77" "  +0.18%  "
# Row 41 - Mantle
Set-Row 41 $null $null "0.781" "  -1.48%  "
# Row 42 - Filecoin
Set-Row 42 $null $null "4.34" "  -2.49%  "
# Row 43 - OKB
Set-Row 43 $null $null "40.39" "  +0.69%  "
# Row 44 - RenderToken
Set-Row 44 $null $null "6.03" "  -3.93%  "
# Row 45 - Hedera
Set-Row 45 $null $null "0.0657" "  -2.12%  "

# Row 46,47,48 rotate: Bittensor -> dogwifhat, InjectiveProtocol -> Bittensor, dogwifhat -> InjectiveProtocol
Set-Row 46 "dogwifhat" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif" "2.31" "  -4.36%  "
Set-Row 47 "Bittensor" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao" "321.24" "  -2.17%  "
Set-Row 48 "InjectiveProtocol" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj" "23.52" "  -5.50%  "

# Row 49 - VeChain
Set-Row 49 $null $null $null "  -2.00%  "
# Row 50 - Stellar
Set-Row 50 $null $null "0.103" "  +2.69%  "
# Row 51 - FirstDigitalUSD
Set-Row 51 $null $null $null "  -0.05%  "
